# Generate Report for Handoff
# Refresh the "Latest Handoff Date/Datetime" values for the
# 4463304d-68d6-4969-b9a8-c4c1c52acd24 file row (row 7) across all sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Per-locale sheets: column E holds "Latest Handoff Datetime"
$zhcn.Range("E7").Value = "2016-03-24 09:43:02"
$dede.Range("E7").Value = "2016-03-24 09:43:08"

# Overview sheet: column D holds "Latest Handoff Date" (max across locales)
$overview.Range("D7").Value = "2016-03-24 09:43:08"
